# Fruta / hortaliza, semanal
# Insert a brand-new weekly record at row 2 (pushing all existing rows down
# by one, old row 29 -> new row 30), matching the latest market report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2 (first data row after the header).
$ws.Rows.Item(2).Insert()

# The Insert() call copies formatting from the header row (bold / centered)
# onto the whole new row. Strip that back to the plain, unstyled look used
# by every other data row, then re-apply the date number format to column D
# only (matching the rest of the "Fecha" column).
$ws.Range("A2:R2").ClearFormats()
$ws.Range("D2").NumberFormat = $ws.Range("D3").NumberFormat

# Populate the new row with the new weekly record.
$ws.Range("A2").Value = 5
$ws.Range("B2").Value = "Macroferia Regional de Talca"
$ws.Range("C2").Value = "Maule"
$ws.Range("D2").Value = 44473
$ws.Range("E2").Value = 7
$ws.Range("F2").Value = 100112026
$ws.Range("G2").Value = "Haba"
$ws.Range("H2").Value = "Sin especificar"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 500
$ws.Range("K2").Value = 9000
$ws.Range("L2").Value = 9000
$ws.Range("M2").Value = 9000
$ws.Range("N2").Value = "$/saco 25 kilos"
$ws.Range("O2").Value = "Región de O'Higgins"
$ws.Range("P2").Value = 360
$ws.Range("Q2").Value = 25
$ws.Range("R2").Value = "Hortaliza"
